# 4.1.1.1a — add the 2022 column (U) mirroring the existing 2021 column (T):
#   - copy T3:T40 formatting into U3:U40 (so borders/number formats/fonts match column T row-by-row)
#   - fill in the 2022 values (numbers, "-" placeholders, and the new header year)
#   - update the view: scroll so column B is the left-most visible column and select V6

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bring the new column's formatting in line with column T (same style per row,
# including the un-valued header/section rows) instead of hand-rolling every
# style id.
$ws.Range("T3:T40").Copy()
$ws.Range("U3:U40").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Header year
$ws.Range("U4").Value = 2022

# "Boys" block (rows 6-22)
$ws.Range("U6").Value = 1456
$ws.Range("U8").Value = 45
$ws.Range("U9").Value = 35
$ws.Range("U10").Value = "-"
$ws.Range("U11").Value = 217
$ws.Range("U12").Value = 22
$ws.Range("U13").Value = 8
$ws.Range("U14").Value = "-"
$ws.Range("U15").Value = "-"
$ws.Range("U16").Value = 57
$ws.Range("U17").Value = "-"
$ws.Range("U18").Value = 5
$ws.Range("U19").Value = "-"
$ws.Range("U20").Value = 46
$ws.Range("U21").Value = 1021
$ws.Range("U22").Value = "-"

# "Girls" block (rows 24-40)
$ws.Range("U24").Value = 1019
$ws.Range("U26").Value = 15
$ws.Range("U27").Value = 30
$ws.Range("U28").Value = 1
$ws.Range("U29").Value = 179
$ws.Range("U30").Value = 16
$ws.Range("U31").Value = 8
$ws.Range("U32").Value = "-"
$ws.Range("U33").Value = "-"
$ws.Range("U34").Value = 46
$ws.Range("U35").Value = "-"
$ws.Range("U36").Value = "-"
$ws.Range("U37").Value = "-"
$ws.Range("U38").Value = 25
$ws.Range("U39").Value = 699
$ws.Range("U40").Value = "-"

# View: scroll so column B is left-most visible, select V6 (matches the saved view state)
$ws.Range("V6").Select()
$excel.ActiveWindow.ScrollColumn = 2
